$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 16
$ws.Cells.Item(2, 2).Value = "Cannot track subject Subject too Large Move Away and retry ."
$ws.Cells.Item(2, 3).Value = "Cannot track subject"
$ws.Cells.Item(2, 4).Value = "0-2"
$ws.Cells.Item(2, 5).Value = "Missing"

$ws.Cells.Item(3, 1).Value = 16
$ws.Cells.Item(3, 2).Value = "Cannot track subject Subject too Large Move Away and retry ."
$ws.Cells.Item(3, 3).Value = "Subject too Large"
$ws.Cells.Item(3, 4).Value = "3-5"
$ws.Cells.Item(3, 5).Value = "Missing"

$ws.Cells.Item(4, 1).Value = 16
$ws.Cells.Item(4, 2).Value = "Cannot track subject Subject too Large Move Away and retry ."
$ws.Cells.Item(4, 3).Value = "Cannot track subject Subject too Large"
$ws.Cells.Item(4, 4).Value = "0-5"
$ws.Cells.Item(4, 5).Formula = "'False"

$ws.Cells.Item(5, 1).Value = 19
$ws.Cells.Item(5, 2).Value = "Switched to gimbal free mode Use the RC to control aircraft yaw ."
$ws.Cells.Item(5, 3).Value = "Use the RC to control aircraft yaw"
$ws.Cells.Item(5, 4).Value = "5-11"
$ws.Cells.Item(5, 5).Value = "Missing"

$ws.Cells.Item(6, 1).Value = 23
$ws.Cells.Item(6, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(6, 3).Value = "Compass Interference"
$ws.Cells.Item(6, 4).Value = "0-1"
$ws.Cells.Item(6, 5).Value = "Missing"

$ws.Cells.Item(7, 1).Value = 23
$ws.Cells.Item(7, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(7, 3).Value = "Temp Max Altitude: nnn"
$ws.Cells.Item(7, 4).Value = "2-5"
$ws.Cells.Item(7, 5).Value = "Missing"

$ws.Cells.Item(8, 1).Value = 23
$ws.Cells.Item(8, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(8, 3).Value = "Compass Interference Temp"
$ws.Cells.Item(8, 4).Value = "0-2"
$ws.Cells.Item(8, 5).Formula = "'False"

$ws.Cells.Item(9, 1).Value = 23
$ws.Cells.Item(9, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(9, 3).Value = "Max Altitude: nnn"
$ws.Cells.Item(9, 4).Value = "3-5"
$ws.Cells.Item(9, 5).Formula = "'False"

$ws.Cells.Item(10, 1).Value = 28
$ws.Cells.Item(10, 2).Value = "Your aircraft has entered a Warning Zone (Class D) Please fly with caution ."
$ws.Cells.Item(10, 3).Value = "Your aircraft has entered a Warning Zone (Class D)"
$ws.Cells.Item(10, 4).Value = "0-8"
$ws.Cells.Item(10, 5).Value = "Missing"

$ws.Cells.Item(11, 1).Value = 28
$ws.Cells.Item(11, 2).Value = "Your aircraft has entered a Warning Zone (Class D) Please fly with caution ."
$ws.Cells.Item(11, 3).Value = "Your aircraft has entered a Warning Zone (Class"
$ws.Cells.Item(11, 4).Value = "0-7"
$ws.Cells.Item(11, 5).Formula = "'False"

$ws.Cells.Item(12, 1).Value = 31
$ws.Cells.Item(12, 2).Value = "Warning: Battery Temperature Below 15°C (59F) Warm battery to above 25°C (77F) before flying ."
$ws.Cells.Item(12, 3).Value = "Warning: Battery Temperature Below 15°C (59F)"
$ws.Cells.Item(12, 4).Value = "0-5"
$ws.Cells.Item(12, 5).Value = "Missing"

$ws.Cells.Item(13, 1).Value = 32
$ws.Cells.Item(13, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(13, 3).Value = "Compass Interference"
$ws.Cells.Item(13, 4).Value = "0-1"
$ws.Cells.Item(13, 5).Value = "Missing"

$ws.Cells.Item(14, 1).Value = 32
$ws.Cells.Item(14, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(14, 3).Value = "Temp Max Altitude 98ft"
$ws.Cells.Item(14, 4).Value = "2-5"
$ws.Cells.Item(14, 5).Value = "Missing"

$ws.Cells.Item(15, 1).Value = 32
$ws.Cells.Item(15, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(15, 3).Value = "Compass Interference Temp"
$ws.Cells.Item(15, 4).Value = "0-2"
$ws.Cells.Item(15, 5).Formula = "'False"

$ws.Cells.Item(16, 1).Value = 32
$ws.Cells.Item(16, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(16, 3).Value = "Max Altitude 98ft"
$ws.Cells.Item(16, 4).Value = "3-5"
$ws.Cells.Item(16, 5).Formula = "'False"

$ws.Cells.Item(17, 1).Value = 33
$ws.Cells.Item(17, 2).Value = "Cannot track subject Subject too Small Get Closer and retry ."
$ws.Cells.Item(17, 3).Value = "Cannot track subject"
$ws.Cells.Item(17, 4).Value = "0-2"
$ws.Cells.Item(17, 5).Value = "Missing"

$ws.Cells.Item(18, 1).Value = 33
$ws.Cells.Item(18, 2).Value = "Cannot track subject Subject too Small Get Closer and retry ."
$ws.Cells.Item(18, 3).Value = "Subject too Small"
$ws.Cells.Item(18, 4).Value = "3-5"
$ws.Cells.Item(18, 5).Value = "Missing"

$ws.Cells.Item(19, 1).Value = 33
$ws.Cells.Item(19, 2).Value = "Cannot track subject Subject too Small Get Closer and retry ."
$ws.Cells.Item(19, 3).Value = "Cannot track subject Subject too Small"
$ws.Cells.Item(19, 4).Value = "0-5"
$ws.Cells.Item(19, 5).Formula = "'False"

$ws.Cells.Item(20, 1).Value = 34
$ws.Cells.Item(20, 2).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) Fly with caution ."
$ws.Cells.Item(20, 3).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant)"
$ws.Cells.Item(20, 4).Value = "0-13"
$ws.Cells.Item(20, 5).Value = "Missing"

$ws.Cells.Item(21, 1).Value = 34
$ws.Cells.Item(21, 2).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) Fly with caution ."
$ws.Cells.Item(21, 3).Value = "GEO: You are in a Warning Zone (Airport"
$ws.Cells.Item(21, 4).Value = "0-7"
$ws.Cells.Item(21, 5).Formula = "'False"

$ws.Cells.Item(22, 1).Value = 34
$ws.Cells.Item(22, 2).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) Fly with caution ."
$ws.Cells.Item(22, 3).Value = "Airports Power Plant)"
$ws.Cells.Item(22, 4).Value = "11-13"
$ws.Cells.Item(22, 5).Formula = "'False"

$ws.Cells.Item(23, 1).Value = 36
$ws.Cells.Item(23, 2).Value = "Incompatible firmware version Go to Profile > Settings to update firmware ."
$ws.Cells.Item(23, 3).Value = "Go to Profile > Settings to update firmware"
$ws.Cells.Item(23, 4).Value = "3-10"
$ws.Cells.Item(23, 5).Value = "Missing"
